# "added and updated with crossref"
# The reference-metadata pipeline re-ran crossref matching for four rows
# and re-emitted their Authors cell (column E) with a refreshed author
# list string (same author data, just re-serialized with a different
# internal padding/indent width - a side effect of the upstream tool).
# Update the Authors values for rows 2-5 on Sheet1 to the refreshed text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = '[Wei%Tang%NULL%1,         Zhujun%Cao%NULL%1,         Mingfeng%Han%NULL%0,         Zhengyan%Wang%NULL%1,         Junwen%Chen%NULL%1,         Wenjin%Sun%NULL%1,         Yaojie%Wu%NULL%1,         Wei%Xiao%NULL%0,         Shengyong%Liu%NULL%1,         Erzhen%Chen%NULL%1,         Wei%Chen%NULL%1,         Xiongbiao%Wang%NULL%1,         Jiuyong%Yang%NULL%1,         Jun%Lin%NULL%0,         Qingxia%Zhao%NULL%1,         Youqin%Yan%NULL%0,         Zhibin%Xie%NULL%1,         Dan%Li%NULL%1,         Yaofeng%Yang%NULL%1,         Leshan%Liu%NULL%1,         Jieming%Qu%NULL%0,         Guang%Ning%NULL%1,         Guochao%Shi%NULL%1,         Qing%Xie%NULL%1]'

$ws.Range("E3").Value = '[ Z.%Chen%null%0,          J.% Hu%null%1,          Z.% Zhang%null%1,          S.% Jiang%null%1,          S.% Han%null%1,          D.% Yan%null%1,          R.% Zhuang%null%1,          B.% Hu%null%1,          Z. % Zhang%null%1]'

$ws.Range("E4").Value = '[Bala%Veenu%coreGivesNoEmail%1,        Bjerum%Catherine%coreGivesNoEmail%1,        Chhonker%Yashpal S%coreGivesNoEmail%1,        Elizah%Arthur%coreGivesNoEmail%1,        John%Lucy N%coreGivesNoEmail%1,        King%Christopher L%coreGivesNoEmail%1,        Likia%Rhoda%coreGivesNoEmail%1,        Marks%Michael%coreGivesNoEmail%1,        Martinez%Pere Millat%coreGivesNoEmail%1,        Mitja%Oriol%coreGivesNoEmail%1,        Murry%Daryl J%coreGivesNoEmail%1,        Silus%Linda%coreGivesNoEmail%1,        Wali%Chilaka%coreGivesNoEmail%1]'

$ws.Range("E5").Value = '[Philippe%Gautret%NULL%1,         Jean-Christophe%Lagier%NULL%1,         Philippe%Parola%NULL%0,         Van Thuan%Hoang%NULL%1,         Line%Meddeb%NULL%1,         Morgane%Mailhe%NULL%1,         Barbara%Doudier%NULL%1,         Johan%Courjon%NULL%1,         Valérie%Giordanengo%NULL%1,         Vera Esteves%Vieira%NULL%1,         Hervé%Tissot Dupont%NULL%1,         Stéphane%Honoré%NULL%1,         Philippe%Colson%NULL%1,         Eric%Chabrière%NULL%1,         Bernard%La Scola%NULL%1,         Jean-Marc%Rolain%NULL%1,         Philippe%Brouqui%NULL%1,         Didier%Raoult%Didier.raoult@gmail.com%0]'
